# Update "want to go" counts (column F) per the 456a3b4 data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3693
$ws.Range("F5").Value = 3693
$ws.Range("F6").Value = 279
$ws.Range("F7").Value = 5221
$ws.Range("F8").Value = 562
$ws.Range("F9").Value = 394
$ws.Range("F10").Value = 216
$ws.Range("F11").Value = 719
$ws.Range("F13").Value = 117
$ws.Range("F15").Value = 720
$ws.Range("F16").Value = 336
$ws.Range("F18").Value = 95
$ws.Range("F19").Value = 163
$ws.Range("F22").Value = 5972
$ws.Range("F26").Value = 6286
$ws.Range("F28").Value = 20
$ws.Range("F29").Value = 3239
$ws.Range("F30").Value = 355
$ws.Range("F31").Value = 731
$ws.Range("F36").Value = 1094
$ws.Range("F39").Value = 3
$ws.Range("F40").Value = 903
$ws.Range("F41").Value = 1071
$ws.Range("F42").Value = 2044

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1140

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1140
$ws.Range("F7").Value = 3693
$ws.Range("F8").Value = 3693
$ws.Range("F9").Value = 279
$ws.Range("F10").Value = 5221
$ws.Range("F11").Value = 562
$ws.Range("F12").Value = 394
$ws.Range("F13").Value = 216
$ws.Range("F14").Value = 719
$ws.Range("F16").Value = 117
$ws.Range("F18").Value = 720
$ws.Range("F19").Value = 336
$ws.Range("F22").Value = 95
$ws.Range("F23").Value = 163
$ws.Range("F26").Value = 5972
$ws.Range("F30").Value = 6286
$ws.Range("F32").Value = 20
$ws.Range("F33").Value = 3239
$ws.Range("F34").Value = 355
$ws.Range("F35").Value = 731
$ws.Range("F41").Value = 1094
$ws.Range("F44").Value = 3
$ws.Range("F45").Value = 903
$ws.Range("F46").Value = 1071
$ws.Range("F48").Value = 2044
